# Update MSME Country Indicators - Luxembourg Summary with refreshed
# (more precise) percentage / density figures.
#
# All of the touched cells hold their numbers as *text* (shared strings,
# not numeric cells), so we force the target range to Text format before
# writing the value and then restore its original style, which keeps the
# cell's type as text (t="s") without introducing a new number-format /
# style index.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, [string]$value) {
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = $origStyle
}

# Row 11 - "Enterprises density (per 1000 people)" (Statistical Institution table)
Set-TextValue $ws.Range("B11") "54.64"
Set-TextValue $ws.Range("C11") "8.18"
Set-TextValue $ws.Range("D11") "62.82"

# Row 33 - "Enterprises density (per 1000 people)" (SME Associations table)
Set-TextValue $ws.Range("B33") "48.97"
Set-TextValue $ws.Range("C33") "7.02"
Set-TextValue $ws.Range("D33") "55.99"

# Row 34 - "Employment (% of total)"
Set-TextValue $ws.Range("B34") "23.11"
Set-TextValue $ws.Range("C34") "42.73"
Set-TextValue $ws.Range("D34") "65.84"

# Row 36 - "Enterprises (% of total)"
Set-TextValue $ws.Range("B36") "87.03"
Set-TextValue $ws.Range("C36") "12.48"
Set-TextValue $ws.Range("D36") "99.51"

# Row 40 - "Value added to the economy (% of total)"
Set-TextValue $ws.Range("B40") "29.88"
Set-TextValue $ws.Range("C40") "39.28"
Set-TextValue $ws.Range("D40") "69.17"
